$wb = $excel.ActiveWorkbook

# --- Step 1: create the new "mRNA" sheet as a duplicate of the existing sheet ---
# The target sheetId for the new sheet is 3 (the workbook already "used" sheetId 1
# for the original sheet). We add a throwaway sheet first to consume sheetId=2,
# then copy the original sheet (which becomes sheetId=3), then delete the
# throwaway sheet, leaving the copy with the desired sheetId and in first position.

# Add a throwaway sheet before the (currently active) original sheet.
$wb.Worksheets.Add() | Out-Null

# The original sheet is now pushed to index 2. Duplicate it (placed right before
# itself, i.e. at index 2, pushing the original to index 3).
$wb.Worksheets.Item(2).Copy($wb.Worksheets.Item(2), $null)

# Remove the throwaway sheet (still at index 1).
$wb.Worksheets.Item(1).Delete() | Out-Null

# Now: index1 = new duplicated sheet (sheetId=3), index2 = original sheet (sheetId=1)
$newSheet = $wb.Worksheets.Item(1)
$origSheet = $wb.Worksheets.Item(2)

# --- Step 2: rename both sheets ---
$newSheet.Name = "Comparison Tables mRNA)"
$origSheet.Name = "Comparison Tables miRNA"

# --- Step 3: clear out the data-value cells in the new mRNA sheet ---
# (the mRNA results were not filled in yet - only the row/column labels remain)
$newSheet.Range("D5:E15").ClearContents()
$newSheet.Range("I5:J8").ClearContents()
$newSheet.Range("N5:O6").ClearContents()
$newSheet.Range("N9:O10").ClearContents()

# --- Step 4: set the active selection on the new sheet ---
$newSheet.Range("N16").Select() | Out-Null
